# Add a new "Oran-Rief Mines" sheet (cloned from the "Xantcha's Crucible"
# template sheet, same layout/styles) and populate it with the Caterium
# Ore -> Caterium Ingots production data, then record the new mine's
# power draw into the "Material Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Clone the "Xantcha's Crucible" sheet as a template and place the
#        new sheet right after it (i.e. as the new last tab). ------------
$template = $wb.Worksheets.Item("Xantcha's Crucible")
$template.Copy($null, $template) | Out-Null
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Oran-Rief Mines"
$template.Range("B2:L36").Select() | Out-Null

# --- 2. Fill in the new sheet's data (order matters so new shared
#        strings land at the same table offsets as the authored file). --
$ws.Range("C2").Value = "Oran-Rief Mines"
$ws.Range("C6").Value = "803.3 MW"
$ws.Range("C7").Value = "-803.3 MW"
$ws.Range("C9").Value = "2.74 Hours"
$ws.Range("C8").Value = "2200 MWh"
$ws.Range("E4").Value = "Caterium Ore"
$ws.Range("J4").Value = "Caterium Ingots"
$ws.Range("F4").Value = 480
$ws.Range("F5").Value = 480
$ws.Range("K4").Value = 240

# --- 3. Record the new mine's output on the "Material Summary" sheet
#        (this sheet tracks ingots/outputs, mirroring row 9's Copper
#        Ingots / Xantcha's Crucible entry). ----------------------------
$summary = $wb.Worksheets.Item("Material Summary")
$summary.Range("C10").Value = "Caterium Ingots"
$summary.Range("D10").Value = 240
$summary.Range("E10").Value = "Oran-Rief Mines"

# --- 4. Make the new sheet the active tab/selection, matching the
#        freshly-added sheet being the one the author was looking at. ---
$ws.Activate() | Out-Null
$ws.Range("C23").Select() | Out-Null
